$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 65 ("「かたちといろ」" entry) - subsequent rows shift up automatically
$ws.Rows.Item(65).Delete()
